# Update cluster-pair LR data for Artn-Gfra1 with new TPM-derived values.
# A new target cluster "ECs" is introduced, expanding the table from 4 to 6 data rows
# (2 sending clusters FAPs/MuSCs x 3 target clusters ECs/FAPs/MuSCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.298902666666667
$ws.Range("H2").Value = 3.896708
$ws.Range("I2").Value = 0.7964389134426562
$ws.Range("J2").Value = 0.7964389134426563
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.048574
$ws.Range("N2").Value = 0.145722
$ws.Range("O2").Value = 0.003596689925699211
$ws.Range("P2").Value = 0.003596689925699211
$ws.Range("Q2").Value = 0.06309289813066667
$ws.Range("R2").Value = 0.567836083176
$ws.Range("S2").Value = 0.002864543816414027
$ws.Range("T2").Value = 0.002864543816414028

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.298902666666667
$ws.Range("H3").Value = 3.896708
$ws.Range("I3").Value = 0.7964389134426562
$ws.Range("J3").Value = 0.7964389134426563
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.922308333333334
$ws.Range("N3").Value = 26.766925
$ws.Range("O3").Value = 0.6606574812962104
$ws.Range("P3").Value = 0.6606574812962104
$ws.Range("Q3").Value = 11.58921008698889
$ws.Range("R3").Value = 104.3028907829
$ws.Range("S3").Value = 0.5261733265613158
$ws.Range("T3").Value = 0.5261733265613158

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.298902666666667
$ws.Range("H4").Value = 3.896708
$ws.Range("I4").Value = 0.7964389134426562
$ws.Range("J4").Value = 0.7964389134426563
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.534313
$ws.Range("N4").Value = 13.602939
$ws.Range("O4").Value = 0.3357458287780905
$ws.Range("P4").Value = 0.3357458287780905
$ws.Range("Q4").Value = 5.889631247201334
$ws.Range("R4").Value = 53.00668122481201
$ws.Range("S4").Value = 0.2674010430649265
$ws.Range("T4").Value = 0.2674010430649265

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3319853333333333
$ws.Range("H5").Value = 0.995956
$ws.Range("I5").Value = 0.2035610865573438
$ws.Range("J5").Value = 0.2035610865573438
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.048574
$ws.Range("N5").Value = 0.145722
$ws.Range("O5").Value = 0.003596689925699211
$ws.Range("P5").Value = 0.003596689925699211
$ws.Range("Q5").Value = 0.01612585558133333
$ws.Range("R5").Value = 0.145132700232
$ws.Range("S5").Value = 0.0007321461092851835
$ws.Range("T5").Value = 0.0007321461092851836

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3319853333333333
$ws.Range("H6").Value = 0.995956
$ws.Range("I6").Value = 0.2035610865573438
$ws.Range("J6").Value = 0.2035610865573438
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.922308333333334
$ws.Range("N6").Value = 26.766925
$ws.Range("O6").Value = 0.6606574812962104
$ws.Range("P6").Value = 0.6606574812962104
$ws.Range("Q6").Value = 2.962075506144444
$ws.Range("R6").Value = 26.6586795553
$ws.Range("S6").Value = 0.1344841547348946
$ws.Range("T6").Value = 0.1344841547348946

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.3319853333333333
$ws.Range("H7").Value = 0.995956
$ws.Range("I7").Value = 0.2035610865573438
$ws.Range("J7").Value = 0.2035610865573438
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.534313
$ws.Range("N7").Value = 13.602939
$ws.Range("O7").Value = 0.3357458287780905
$ws.Range("P7").Value = 0.3357458287780905
$ws.Range("Q7").Value = 1.505325412742667
$ws.Range("R7").Value = 13.547928714684
$ws.Range("S7").Value = 0.068344785713164
$ws.Range("T7").Value = 0.06834478571316402

